# Updated 2D training schedules, no break screen
# Adds a new "break_on_off" column (L) to the trial schedule sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header for the new column
$ws.Cells.Item(1, 12).Value = "break_on_off"

# Rows that mark a break (on/off = 1); every other data row is 0
$breakRows = @(19, 37, 54)

for ($r = 2; $r -le 73; $r++) {
    if ($breakRows -contains $r) {
        $ws.Cells.Item($r, 12).Value = 1
    } else {
        $ws.Cells.Item($r, 12).Value = 0
    }
}

# Match the saved selection: whole new column selected, active cell at the top
$ws.Range("L1:L73").Select() | Out-Null
